# Change the A1 heading cell from "Heading 1" to "Format".
# All the other cells keep their existing (string) values; the shared-string
# table re-indexing that results from removing the old "Heading 1" entry and
# appending the new "Format" entry happens automatically when the workbook
# is saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Format"
